$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new blank column before the "Late" column ---
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = 10.33

# Repayment Schedule becomes the active / selected sheet & cell
$wsSchedule.Activate()
$wsSchedule.Range("J19").Select()
